$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "73×54=3942" "55×14=770"
Replace-Text "23×29=667" "40×16=640"
Replace-Text "98×85=8330" "32×65=2080"
Replace-Text "56×98=5488" "89×29=2581"
Replace-Text "26×26=676" "52×99=5148"
Replace-Text "87×85=7395" "31×22=682"
Replace-Text "55×35=1925" "67×26=1742"
Replace-Text "76×33=2508" "97×71=6887"
Replace-Text "65×99=6435" "46×16=736"
Replace-Text "51×73=3723" "99×38=3762"
Replace-Text "16×28=448" "79×16=1264"
Replace-Text "25×77=1925" "35×30=1050"
Replace-Text "57×40=2280" "72×32=2304"
Replace-Text "84×34=2856" "59×95=5605"
Replace-Text "78×67=5226" "77×16=1232"
Replace-Text "25×70=1750" "58×49=2842"
Replace-Text "60×57=3420" "99×23=2277"
Replace-Text "36×26=936" "67×40=2680"
Replace-Text "29×11=319" "92×92=8464"
Replace-Text "57×24=1368" "28×71=1988"
Replace-Text "61×21=1281" "84×21=1764"
Replace-Text "42×86=3612" "49×34=1666"
Replace-Text "14×63=882" "66×30=1980"
Replace-Text "95×79=7505" "72×75=5400"
Replace-Text "95×83=7885" "57×54=3078"
